function Set-TextCell($cell, [string]$val) {
    $cell.Value = "'" + $val
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# A sheet that still carries the original header/index-column style (s=2 in
# styles.xml) verbatim -- used below to stamp the same style onto newly
# created cells via Copy + PasteSpecial(xlPasteFormats).
$srcFmt = $wb.Worksheets.Item("2021-Q4")

# ---- Step 1: restructure sheets ----
# Remove the existing "总计" sheet; it will be re-added at the end so that it
# (re)receives the next sheetId in sequence, matching the target where the
# new "2022-Q1" sheet takes sheetId 4 / sheet4.xml and "总计" becomes
# sheetId 5 / sheet5.xml.
$zj = $wb.Worksheets.Item("总计")
$zj.Delete()

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs = $wb.Worksheets.Add()
$newWs.Name = "2022-Q1"
$newWs.Move([System.Reflection.Missing]::Value, $lastSheet)

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$zjWs = $wb.Worksheets.Add()
$zjWs.Name = "总计"
$zjWs.Move([System.Reflection.Missing]::Value, $lastSheet2)


# ---- Step 3: populate "2022-Q1" sheet (new sheet4.xml) ----
# ---- header row for 2022-Q1 sheet (B1:H1) ----
$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# ---- data rows 2..25 ----
$newWs.Cells.Item(2,1).Value = 0
Set-TextCell $newWs.Cells.Item(2,2) "007130"
$newWs.Cells.Item(2,3).Value = "中庚小盘价值股票"
Set-TextCell $newWs.Cells.Item(2,4) "40.99"
Set-TextCell $newWs.Cells.Item(2,5) "93.10"
Set-TextCell $newWs.Cells.Item(2,6) "7.51"
Set-TextCell $newWs.Cells.Item(2,7) "3.0783"
$newWs.Cells.Item(2,8).Value = 1
$newWs.Cells.Item(3,1).Value = 1
Set-TextCell $newWs.Cells.Item(3,2) "006551"
$newWs.Cells.Item(3,3).Value = "中庚价值领航混合"
Set-TextCell $newWs.Cells.Item(3,4) "36.49"
Set-TextCell $newWs.Cells.Item(3,5) "94.18"
Set-TextCell $newWs.Cells.Item(3,6) "4.14"
Set-TextCell $newWs.Cells.Item(3,7) "1.5107"
$newWs.Cells.Item(3,8).Value = 8
$newWs.Cells.Item(4,1).Value = 2
Set-TextCell $newWs.Cells.Item(4,2) "011174"
$newWs.Cells.Item(4,3).Value = "中庚价值品质一年持有期混合"
Set-TextCell $newWs.Cells.Item(4,4) "63.01"
Set-TextCell $newWs.Cells.Item(4,5) "93.47"
Set-TextCell $newWs.Cells.Item(4,6) "2.28"
Set-TextCell $newWs.Cells.Item(4,7) "1.4366"
$newWs.Cells.Item(4,8).Value = 8
$newWs.Cells.Item(5,1).Value = 3
Set-TextCell $newWs.Cells.Item(5,2) "007497"
$newWs.Cells.Item(5,3).Value = "中庚价值灵动灵活配置混合"
Set-TextCell $newWs.Cells.Item(5,4) "24.35"
Set-TextCell $newWs.Cells.Item(5,5) "89.42"
Set-TextCell $newWs.Cells.Item(5,6) "5.30"
Set-TextCell $newWs.Cells.Item(5,7) "1.2906"
$newWs.Cells.Item(5,8).Value = 1
$newWs.Cells.Item(6,1).Value = 4
Set-TextCell $newWs.Cells.Item(6,2) "005392"
$newWs.Cells.Item(6,3).Value = "长信价值蓝筹两年定期开放灵活配置混合A"
Set-TextCell $newWs.Cells.Item(6,4) "5.17"
Set-TextCell $newWs.Cells.Item(6,5) "93.17"
Set-TextCell $newWs.Cells.Item(6,6) "8.03"
Set-TextCell $newWs.Cells.Item(6,7) "0.4152"
$newWs.Cells.Item(6,8).Value = 3
$newWs.Cells.Item(7,1).Value = 5
Set-TextCell $newWs.Cells.Item(7,2) "009911"
$newWs.Cells.Item(7,3).Value = "长信价值蓝筹两年定期开放灵活配置混合C"
Set-TextCell $newWs.Cells.Item(7,4) "4.99"
Set-TextCell $newWs.Cells.Item(7,5) "93.17"
Set-TextCell $newWs.Cells.Item(7,6) "8.03"
Set-TextCell $newWs.Cells.Item(7,7) "0.4007"
$newWs.Cells.Item(7,8).Value = 3
$newWs.Cells.Item(8,1).Value = 6
Set-TextCell $newWs.Cells.Item(8,2) "012426"
$newWs.Cells.Item(8,3).Value = "南方价值臻选混合型证券投资基金A"
Set-TextCell $newWs.Cells.Item(8,4) "5.60"
Set-TextCell $newWs.Cells.Item(8,5) "57.78"
Set-TextCell $newWs.Cells.Item(8,6) "1.04"
Set-TextCell $newWs.Cells.Item(8,7) "0.0582"
$newWs.Cells.Item(8,8).Value = 4
$newWs.Cells.Item(9,1).Value = 7
Set-TextCell $newWs.Cells.Item(9,2) "160324"
$newWs.Cells.Item(9,3).Value = "华夏磐晟灵活配置混合（LOF）"
Set-TextCell $newWs.Cells.Item(9,4) "1.07"
Set-TextCell $newWs.Cells.Item(9,5) "78.29"
Set-TextCell $newWs.Cells.Item(9,6) "4.73"
Set-TextCell $newWs.Cells.Item(9,7) "0.0506"
$newWs.Cells.Item(9,8).Value = 6
$newWs.Cells.Item(10,1).Value = 8
Set-TextCell $newWs.Cells.Item(10,2) "000524"
$newWs.Cells.Item(10,3).Value = "上投摩根民生需求股票"
Set-TextCell $newWs.Cells.Item(10,4) "1.57"
Set-TextCell $newWs.Cells.Item(10,5) "83.43"
Set-TextCell $newWs.Cells.Item(10,6) "3.21"
Set-TextCell $newWs.Cells.Item(10,7) "0.0504"
$newWs.Cells.Item(10,8).Value = 5
$newWs.Cells.Item(11,1).Value = 9
Set-TextCell $newWs.Cells.Item(11,2) "001181"
$newWs.Cells.Item(11,3).Value = "南方改革机遇灵活配置混合"
Set-TextCell $newWs.Cells.Item(11,4) "4.06"
Set-TextCell $newWs.Cells.Item(11,5) "56.27"
Set-TextCell $newWs.Cells.Item(11,6) "1.09"
Set-TextCell $newWs.Cells.Item(11,7) "0.0443"
$newWs.Cells.Item(11,8).Value = 5
$newWs.Cells.Item(12,1).Value = 10
Set-TextCell $newWs.Cells.Item(12,2) "012244"
$newWs.Cells.Item(12,3).Value = "广发金融地产精选股票型发起式证券投资基金A"
Set-TextCell $newWs.Cells.Item(12,4) "0.98"
Set-TextCell $newWs.Cells.Item(12,5) "86.81"
Set-TextCell $newWs.Cells.Item(12,6) "3.65"
Set-TextCell $newWs.Cells.Item(12,7) "0.0358"
$newWs.Cells.Item(12,8).Value = 9
$newWs.Cells.Item(13,1).Value = 11
Set-TextCell $newWs.Cells.Item(13,2) "001162"
$newWs.Cells.Item(13,3).Value = "前海开源优势蓝筹股票A"
Set-TextCell $newWs.Cells.Item(13,4) "0.62"
Set-TextCell $newWs.Cells.Item(13,5) "92.32"
Set-TextCell $newWs.Cells.Item(13,6) "5.10"
Set-TextCell $newWs.Cells.Item(13,7) "0.0316"
$newWs.Cells.Item(13,8).Value = 5
$newWs.Cells.Item(14,1).Value = 12
Set-TextCell $newWs.Cells.Item(14,2) "001536"
$newWs.Cells.Item(14,3).Value = "南方君选灵活配置混合"
Set-TextCell $newWs.Cells.Item(14,4) "3.82"
Set-TextCell $newWs.Cells.Item(14,5) "48.02"
Set-TextCell $newWs.Cells.Item(14,6) "0.76"
Set-TextCell $newWs.Cells.Item(14,7) "0.0290"
$newWs.Cells.Item(14,8).Value = 9
$newWs.Cells.Item(15,1).Value = 13
Set-TextCell $newWs.Cells.Item(15,2) "000433"
$newWs.Cells.Item(15,3).Value = "安信鑫发优选混合"
Set-TextCell $newWs.Cells.Item(15,4) "1.23"
Set-TextCell $newWs.Cells.Item(15,5) "67.20"
Set-TextCell $newWs.Cells.Item(15,6) "2.15"
Set-TextCell $newWs.Cells.Item(15,7) "0.0264"
$newWs.Cells.Item(15,8).Value = 10
$newWs.Cells.Item(16,1).Value = 14
Set-TextCell $newWs.Cells.Item(16,2) "202213"
$newWs.Cells.Item(16,3).Value = "南方核心竞争混合"
Set-TextCell $newWs.Cells.Item(16,4) "2.80"
Set-TextCell $newWs.Cells.Item(16,5) "45.89"
Set-TextCell $newWs.Cells.Item(16,6) "0.94"
Set-TextCell $newWs.Cells.Item(16,7) "0.0263"
$newWs.Cells.Item(16,8).Value = 5
$newWs.Cells.Item(17,1).Value = 15
Set-TextCell $newWs.Cells.Item(17,2) "012245"
$newWs.Cells.Item(17,3).Value = "广发金融地产精选股票型发起式证券投资基金C"
Set-TextCell $newWs.Cells.Item(17,4) "0.63"
Set-TextCell $newWs.Cells.Item(17,5) "86.81"
Set-TextCell $newWs.Cells.Item(17,6) "3.65"
Set-TextCell $newWs.Cells.Item(17,7) "0.0230"
$newWs.Cells.Item(17,8).Value = 9
$newWs.Cells.Item(18,1).Value = 16
Set-TextCell $newWs.Cells.Item(18,2) "001482"
$newWs.Cells.Item(18,3).Value = "上投摩根新兴服务股票"
Set-TextCell $newWs.Cells.Item(18,4) "0.42"
Set-TextCell $newWs.Cells.Item(18,5) "86.98"
Set-TextCell $newWs.Cells.Item(18,6) "2.75"
Set-TextCell $newWs.Cells.Item(18,7) "0.0116"
$newWs.Cells.Item(18,8).Value = 9
$newWs.Cells.Item(19,1).Value = 17
Set-TextCell $newWs.Cells.Item(19,2) "620007"
$newWs.Cells.Item(19,3).Value = "金元顺安优质精选灵活配置混合A"
Set-TextCell $newWs.Cells.Item(19,4) "0.75"
Set-TextCell $newWs.Cells.Item(19,5) "39.31"
Set-TextCell $newWs.Cells.Item(19,6) "1.45"
Set-TextCell $newWs.Cells.Item(19,7) "0.0109"
$newWs.Cells.Item(19,8).Value = 8
$newWs.Cells.Item(20,1).Value = 18
Set-TextCell $newWs.Cells.Item(20,2) "001375"
$newWs.Cells.Item(20,3).Value = "金元顺安优质精选灵活配置混合C"
Set-TextCell $newWs.Cells.Item(20,4) "0.69"
Set-TextCell $newWs.Cells.Item(20,5) "39.31"
Set-TextCell $newWs.Cells.Item(20,6) "1.45"
Set-TextCell $newWs.Cells.Item(20,7) "0.0100"
$newWs.Cells.Item(20,8).Value = 8
$newWs.Cells.Item(21,1).Value = 19
Set-TextCell $newWs.Cells.Item(21,2) "001638"
$newWs.Cells.Item(21,3).Value = "前海开源优势蓝筹股票C"
Set-TextCell $newWs.Cells.Item(21,4) "0.09"
Set-TextCell $newWs.Cells.Item(21,5) "92.32"
Set-TextCell $newWs.Cells.Item(21,6) "5.10"
Set-TextCell $newWs.Cells.Item(21,7) "0.0046"
$newWs.Cells.Item(21,8).Value = 5
$newWs.Cells.Item(22,1).Value = 20
Set-TextCell $newWs.Cells.Item(22,2) "410009"
$newWs.Cells.Item(22,3).Value = "华富量子生命力混合"
Set-TextCell $newWs.Cells.Item(22,4) "0.10"
Set-TextCell $newWs.Cells.Item(22,5) "91.31"
Set-TextCell $newWs.Cells.Item(22,6) "3.63"
Set-TextCell $newWs.Cells.Item(22,7) "0.0036"
$newWs.Cells.Item(22,8).Value = 5
$newWs.Cells.Item(23,1).Value = 21
Set-TextCell $newWs.Cells.Item(23,2) "012427"
$newWs.Cells.Item(23,3).Value = "南方价值臻选混合型证券投资基金C"
Set-TextCell $newWs.Cells.Item(23,4) "0.22"
Set-TextCell $newWs.Cells.Item(23,5) "57.78"
Set-TextCell $newWs.Cells.Item(23,6) "1.04"
Set-TextCell $newWs.Cells.Item(23,7) "0.0023"
$newWs.Cells.Item(23,8).Value = 4
$newWs.Cells.Item(24,1).Value = 22
Set-TextCell $newWs.Cells.Item(24,2) "006601"
$newWs.Cells.Item(24,3).Value = "国融融泰灵活配置混合A"
Set-TextCell $newWs.Cells.Item(24,4) "0.04"
Set-TextCell $newWs.Cells.Item(24,5) "47.44"
Set-TextCell $newWs.Cells.Item(24,6) "3.59"
Set-TextCell $newWs.Cells.Item(24,7) "0.0014"
$newWs.Cells.Item(24,8).Value = 5
$newWs.Cells.Item(25,1).Value = 23
Set-TextCell $newWs.Cells.Item(25,2) "006602"
$newWs.Cells.Item(25,3).Value = "国融融泰灵活配置混合C"
Set-TextCell $newWs.Cells.Item(25,4) "0.01"
Set-TextCell $newWs.Cells.Item(25,5) "47.44"
Set-TextCell $newWs.Cells.Item(25,6) "3.59"
Set-TextCell $newWs.Cells.Item(25,7) "0.0004"
$newWs.Cells.Item(25,8).Value = 5

# Apply header style (s=2) and index-column style (s=2) to match the
# formatting used throughout the workbook for this table shape.
$srcFmt.Range("B1").Copy()
$newWs.Range("B1:H1").PasteSpecial(-4122)
$srcFmt.Range("A2").Copy()
$newWs.Range("A2:A25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Step 2: populate "总计" sheet (new sheet5.xml) ----
# Header row, same as before.
$zjWs.Range("B1").Value = "日期"
$zjWs.Range("C1").Value = "持有数量(只)"
$zjWs.Range("D1").Value = "持有市值(亿元)"

# Data: new 2022-Q1 summary row on top, followed by the previously existing
# quarters shifted down by one row.
$zjWs.Cells.Item(2,1).Value = 0
$zjWs.Cells.Item(2,2).Value = "2022-Q1"
$zjWs.Cells.Item(2,3).Value = 24
$zjWs.Cells.Item(2,4).Value = 8.550000000000001

$zjWs.Cells.Item(3,1).Value = 1
$zjWs.Cells.Item(3,2).Value = "2021-Q4"
$zjWs.Cells.Item(3,3).Value = 9
$zjWs.Cells.Item(3,4).Value = 8.44

$zjWs.Cells.Item(4,1).Value = 2
$zjWs.Cells.Item(4,2).Value = "2021-Q3"
$zjWs.Cells.Item(4,3).Value = 8
$zjWs.Cells.Item(4,4).Value = 6.96

$zjWs.Cells.Item(5,1).Value = 3
$zjWs.Cells.Item(5,2).Value = "2021-Q2"
$zjWs.Cells.Item(5,3).Value = 3
$zjWs.Cells.Item(5,4).Value = 0.34

# Re-apply the header/index-column style (s=2 in the original workbook) by
# copying formats from a sheet that still carries it verbatim.
$srcFmt.Range("B1").Copy()
$zjWs.Range("B1:D1").PasteSpecial(-4122)
$srcFmt.Range("A2").Copy()
$zjWs.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
